# Apply hybrid bold + color highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet /
# impact paragraphs, matching the target diff.

$d = $word.ActiveDocument
$HighlightColor = 5258796   # RGB(0x2C,0x3E,0x50) == "2C3E50"
$PlusMinus = [char]0x00B1

# Each entry: a unique substring used to locate the target paragraph,
# followed by the ordered list of substrings (within that paragraph)
# that should become bold + colored.
$edits = @(
    @{
        Match = "developed geospatial machine learning algorithms improving demographic"
        Terms = @("23%", "64%")
    },
    @{
        Match = "reducing polling error margins"
        Terms = @("87%", "71%", ($PlusMinus + "4.2%"), ($PlusMinus + "2.1%"))
    },
    @{
        Match = "Wrote RFP and analyzed bids from"
        Terms = @("1,200")
    },
    @{
        Match = "Polling Consortium Database"
        Terms = @("`$400M", "`$1B")
    },
    @{
        Match = "Algorithm reduced mapping costs"
        Terms = @("73.5%", "`$4.7M")
    },
    @{
        Match = "industry standard of 71%"
        Exclude = "reducing polling error margins"
        Terms = @("87%", "71%")
    }
)

foreach ($edit in $edits) {
    $target = $null
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $text = $p.Range.Text
        if ($text.Contains($edit.Match)) {
            if ($edit.Exclude -and $text.Contains($edit.Exclude)) {
                continue
            }
            $target = $p
        }
    }

    if ($target -eq $null) {
        Write-Output "NOT FOUND: $($edit.Match)"
        continue
    }

    $pEnd = $target.Range.End
    $cursor = $target.Range.Start

    foreach ($term in $edit.Terms) {
        $rng = $d.Range($cursor, $pEnd)
        $found = $rng.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Bold = 1
            $rng.Font.Color = $HighlightColor
            $cursor = $rng.End
        } else {
            Write-Output "TERM NOT FOUND: $term (in paragraph matching '$($edit.Match)')"
        }
    }
}

Write-Output "Highlighting complete"
